$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-06 20:28:13", 0.003),
    @("2023-12-06 20:28:51", 0.002),
    @("2023-12-06 20:29:10", 0.0008)
)

$startRow = 43
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
